$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each target cell is written with a leading apostrophe so Excel
# stores numeric-looking strings (e.g. "245.27", "1.00") as literal
# text instead of silently converting them to numbers - matching the
# original workbook, where every data cell is an inline/shared string.
# ClearFormats() then strips the transient "quote prefix" cell style
# that the apostrophe leaves behind, so no cell style index changes.

$ws.Range("D2").Value = '''36.870.62'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '''  -0.68%  '
$ws.Range("E2").ClearFormats()

$ws.Range("D3").Value = '''2.083.68'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '''  +1.68%  '
$ws.Range("E3").ClearFormats()

$ws.Range("E4").Value = '''  -0.15%  '
$ws.Range("E4").ClearFormats()

$ws.Range("D5").Value = '''245.27'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '''  -1.10%  '
$ws.Range("E5").ClearFormats()

$ws.Range("D6").Value = '''0.653'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '''  -1.49%  '
$ws.Range("E6").ClearFormats()

$ws.Range("D8").Value = '''55.58'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '''  -4.44%  '
$ws.Range("E8").ClearFormats()

$ws.Range("D9").Value = '''59.88'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '''  +0.05%  '
$ws.Range("E9").ClearFormats()

$ws.Range("D10").Value = '''0.368'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '''  -3.38%  '
$ws.Range("E10").ClearFormats()

$ws.Range("D11").Value = '''0.0764'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '''  -1.54%  '
$ws.Range("E11").ClearFormats()

$ws.Range("E12").Value = '''  +1.29%  '
$ws.Range("E12").ClearFormats()

$ws.Range("D13").Value = '''15.07'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '''  -4.90%  '
$ws.Range("E13").ClearFormats()

$ws.Range("D14").Value = '''0.889'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '''  +5.59%  '
$ws.Range("E14").ClearFormats()

$ws.Range("D15").Value = '''2.391.98'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '''  +1.78%  '
$ws.Range("E15").ClearFormats()

$ws.Range("D16").Value = '''2.237.39'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '''  +9.09%  '
$ws.Range("E16").ClearFormats()

$ws.Range("D17").Value = '''5.52'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '''  -3.31%  '
$ws.Range("E17").ClearFormats()

$ws.Range("D18").Value = '''36.830.15'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '''  -0.98%  '
$ws.Range("E18").ClearFormats()

$ws.Range("D19").Value = '''17.40'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '''  -3.55%  '
$ws.Range("E19").ClearFormats()

$ws.Range("D20").Value = '''73.08'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '''  -2.35%  '
$ws.Range("E20").ClearFormats()

$ws.Range("D21").Value = '''0.0₃0882'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '''  -1.24%  '
$ws.Range("E21").ClearFormats()

$ws.Range("D22").Value = '''5.47'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '''  +2.71%  '
$ws.Range("E22").ClearFormats()

$ws.Range("D23").Value = '''237.99'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '''  +0.52%  '
$ws.Range("E23").ClearFormats()

$ws.Range("D24").Value = '''1.00'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '''  +0.09%  '
$ws.Range("E24").ClearFormats()

$ws.Range("D25").Value = '''2.41'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '''  -1.99%  '
$ws.Range("E25").ClearFormats()

$ws.Range("E26").Value = '''  +4.69%  '
$ws.Range("E26").ClearFormats()

$ws.Range("D27").Value = '''2.17'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '''  -0.06%  '
$ws.Range("E27").ClearFormats()

$ws.Range("D28").Value = '''168.10'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '''  -0.92%  '
$ws.Range("E28").ClearFormats()

$ws.Range("D29").Value = '''20.63'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '''  +2.99%  '
$ws.Range("E29").ClearFormats()

$ws.Range("D30").Value = '''5.35'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '''  +11.54%  '
$ws.Range("E30").ClearFormats()

$ws.Range("E31").Value = '''  -0.58%  '
$ws.Range("E31").ClearFormats()

$ws.Range("D32").Value = '''1.20'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '''  +8.37%  '
$ws.Range("E32").ClearFormats()

$ws.Range("E33").Value = '''  +5.42%  '
$ws.Range("E33").ClearFormats()

$ws.Range("D34").Value = '''0.0611'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '''  -0.93%  '
$ws.Range("E34").ClearFormats()

$ws.Range("D35").Value = '''2.38'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '''  +5.62%  '
$ws.Range("E35").ClearFormats()

$ws.Range("E36").Value = '''  +0.13%  '
$ws.Range("E36").ClearFormats()

$ws.Range("E37").Value = '''  +4.04%  '
$ws.Range("E37").ClearFormats()

$ws.Range("E38").Value = '''  -5.78%  '
$ws.Range("E38").ClearFormats()

$ws.Range("E39").Value = '''  -3.76%  '
$ws.Range("E39").ClearFormats()

$ws.Range("E40").Value = '''  +2.02%  '
$ws.Range("E40").ClearFormats()

$ws.Range("E41").Value = '''  -0.14%  '
$ws.Range("E41").ClearFormats()

$ws.Range("D42").Value = '''4.86'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '''  -6.57%  '
$ws.Range("E42").ClearFormats()

$ws.Range("D43").Value = '''0.0951'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '''  -5.55%  '
$ws.Range("E43").ClearFormats()

$ws.Range("D44").Value = '''96.93'
$ws.Range("D44").ClearFormats()

$ws.Range("D45").Value = '''2.84'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '''  -12.66%  '
$ws.Range("E45").ClearFormats()

$ws.Range("E46").Value = '''  -6.02%  '
$ws.Range("E46").ClearFormats()

$ws.Range("D47").Value = '''1.352.45'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '''  +6.03%  '
$ws.Range("E47").ClearFormats()

$ws.Range("D48").Value = '''2.44'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '''  -0.18%  '
$ws.Range("E48").ClearFormats()

$ws.Range("B49").Value = '''FraxShare'
$ws.Range("B49").ClearFormats()
$ws.Range("C49").Value = '''https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("C49").ClearFormats()
$ws.Range("D49").Value = '''7.16'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '''  +5.29%  '
$ws.Range("E49").ClearFormats()

$ws.Range("B50").Value = '''MXToken'
$ws.Range("B50").ClearFormats()
$ws.Range("C50").Value = '''https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("C50").ClearFormats()
$ws.Range("D50").Value = '''2.91'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '''  +1.99%  '
$ws.Range("E50").ClearFormats()

$ws.Range("D51").Value = '''2.276.79'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '''  +1.89%  '
$ws.Range("E51").ClearFormats()
